$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomePage")

# Fill in row 7 (new locator: btn_search / Name / submit_search)
$ws.Range("A7").Value = "btn_search"
$ws.Range("B7").Value = "Name"
$ws.Range("C7").Value = "submit_search"

# Fill in row 8 (new locator: lbl_search_result_locator / XPath / product list xpath)
# Write C8 before A8 so new shared strings are appended in the same order
# as the source workbook (xpath text before the label text).
$ws.Range("C8").Value = "//ul[@class='product_list grid row']//a[contains(text(),'Faded Short Sleeve T-shirts')]"
$ws.Range("B8").Value = "XPath"
$ws.Range("A8").Value = "lbl_search_result_locator"

# Move the active cell selection on the HomePage sheet to A8
$ws.Activate()
$ws.Range("A8").Select()
